$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: fill in blank B4 and E4 with "-1"
$ws.Range("B4").Value = "-1"
$ws.Range("E4").Value = "-1"

# Row 5: fix C5 from "abc" to "125"
$ws.Range("C5").Value = "125"

# Row 6: update A6 label text
$ws.Range("A6").Value = "Coefficient of thermal expansion (microstrain/K)"
